# Update with harris (12/1) and ifop (11/28) polls
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the sample sizes (n, column I) for the existing ifop poll rows 132-134
$ws.Cells.Item(132, 9).Value = 1017
$ws.Cells.Item(133, 9).Value = 1007
$ws.Cells.Item(134, 9).Value = 1004

# Column map:
# A id, B year, C week, D month, E day, F firm, G collectmode, H unsure,
# I n, J c_poutou, K c_arthaud, L c_melenchon, M c_roussel, N c_montebourg,
# O c_jadot, P c_hidalgo, Q c_macron, R c_pecresse, S c_barnier,
# T c_bertrand, U c_lassalle, V c_daignant, W c_lepen, X c_zemmour,
# Y c_asselineau, Z c_poisson, AA c_philippot, AB c_lagarde

# New row 135 - harris poll
$r = 135
$ws.Cells.Item($r, 1).Value = 41
$ws.Cells.Item($r, 2).Value = 2021
$ws.Cells.Item($r, 3).Value = 14
$ws.Cells.Item($r, 4).Value = 11
$ws.Cells.Item($r, 5).Value = 28
$ws.Cells.Item($r, 6).Value = "harris"
$ws.Cells.Item($r, 7).Value = "online"
$ws.Cells.Item($r, 8).Value = "included"
$ws.Cells.Item($r, 9).Value = 1801
$ws.Cells.Item($r, 10).Value = 1
$ws.Cells.Item($r, 11).Value = 1
$ws.Cells.Item($r, 12).Value = 10
$ws.Cells.Item($r, 13).Value = 2
$ws.Cells.Item($r, 14).Value = 2
$ws.Cells.Item($r, 15).Value = 7
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 23
$ws.Cells.Item($r, 20).Value = 14
$ws.Cells.Item($r, 21).Value = "T_0.5"
$ws.Cells.Item($r, 22).Value = 2
$ws.Cells.Item($r, 23).Value = 19
$ws.Cells.Item($r, 24).Value = 13
$ws.Cells.Item($r, 25).Value = "T_0.5"
$ws.Cells.Item($r, 25).Font.Color = $ws.Cells.Item(41, 25).Font.Color
$ws.Cells.Item($r, 27).Value = 1

# New row 136 - harris poll
$r = 136
$ws.Cells.Item($r, 1).Value = 41
$ws.Cells.Item($r, 2).Value = 2021
$ws.Cells.Item($r, 3).Value = 14
$ws.Cells.Item($r, 4).Value = 11
$ws.Cells.Item($r, 5).Value = 28
$ws.Cells.Item($r, 6).Value = "harris"
$ws.Cells.Item($r, 7).Value = "online"
$ws.Cells.Item($r, 8).Value = "included"
$ws.Cells.Item($r, 9).Value = 1781
$ws.Cells.Item($r, 10).Value = 1
$ws.Cells.Item($r, 11).Value = 1
$ws.Cells.Item($r, 12).Value = 10
$ws.Cells.Item($r, 13).Value = 2
$ws.Cells.Item($r, 14).Value = 2
$ws.Cells.Item($r, 15).Value = 8
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 24
$ws.Cells.Item($r, 18).Value = 11
$ws.Cells.Item($r, 21).Value = "T_0.5"
$ws.Cells.Item($r, 22).Value = 2
$ws.Cells.Item($r, 23).Value = 20
$ws.Cells.Item($r, 24).Value = 13
$ws.Cells.Item($r, 25).Value = "T_0.5"
$ws.Cells.Item($r, 25).Font.Color = $ws.Cells.Item(41, 25).Font.Color
$ws.Cells.Item($r, 27).Value = 1

# New row 137 - harris poll
$r = 137
$ws.Cells.Item($r, 1).Value = 41
$ws.Cells.Item($r, 2).Value = 2021
$ws.Cells.Item($r, 3).Value = 14
$ws.Cells.Item($r, 4).Value = 11
$ws.Cells.Item($r, 5).Value = 28
$ws.Cells.Item($r, 6).Value = "harris"
$ws.Cells.Item($r, 7).Value = "online"
$ws.Cells.Item($r, 8).Value = "included"
$ws.Cells.Item($r, 9).Value = 1781
$ws.Cells.Item($r, 10).Value = 1
$ws.Cells.Item($r, 11).Value = 1
$ws.Cells.Item($r, 12).Value = 10
$ws.Cells.Item($r, 13).Value = 2
$ws.Cells.Item($r, 14).Value = 2
$ws.Cells.Item($r, 15).Value = 8
$ws.Cells.Item($r, 16).Value = 5
$ws.Cells.Item($r, 17).Value = 24
$ws.Cells.Item($r, 19).Value = 10
$ws.Cells.Item($r, 21).Value = "T_0.5"
$ws.Cells.Item($r, 22).Value = 2
$ws.Cells.Item($r, 23).Value = 20
$ws.Cells.Item($r, 24).Value = 13
$ws.Cells.Item($r, 25).Value = "T_0.5"
$ws.Cells.Item($r, 25).Font.Color = $ws.Cells.Item(41, 25).Font.Color
$ws.Cells.Item($r, 27).Value = 2

# Update the frozen pane / selection to match the newly added rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 130
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I137").Select()
